$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row data: row index, Column A (consumption MW), Column B (timestamp serial)
$data = @(
    @(2, 6064, 46021.95833333334),
    @(3, 5998, 46021.96875),
    @(4, 5876, 46021.97916666666),
    @(5, 5837, 46021.98958333334),
    @(6, 5819, 46022),
    @(7, 5748, 46022.01041666666),
    @(8, 5704, 46022.02083333334),
    @(9, 5663, 46022.03125),
    @(10, 5569, 46022.04166666666),
    @(11, 5523, 46022.05208333334),
    @(12, 5487, 46022.0625),
    @(13, 5459, 46022.07291666666),
    @(14, 5418, 46022.08333333334),
    @(15, 5409, 46022.09375),
    @(16, 5387, 46022.10416666666),
    @(17, 5368, 46022.11458333334),
    @(18, 5355, 46022.125),
    @(19, 5403, 46022.13541666666),
    @(20, 5434, 46022.14583333334),
    @(21, 5425, 46022.15625),
    @(22, 5478, 46022.16666666666),
    @(23, 5516, 46022.17708333334),
    @(24, 5553, 46022.1875),
    @(25, 5606, 46022.19791666666),
    @(26, 5723, 46022.20833333334),
    @(27, 5804, 46022.21875),
    @(28, 5865, 46022.22916666666),
    @(29, 5920, 46022.23958333334),
    @(30, 5950, 46022.25),
    @(31, 6051, 46022.26041666666),
    @(32, 6146, 46022.27083333334),
    @(33, 6243, 46022.28125),
    @(34, 6346, 46022.29166666666),
    @(35, 6419, 46022.30208333334),
    @(36, 6477, 46022.3125),
    @(37, 6497, 46022.32291666666),
    @(38, 6519, 46022.33333333334),
    @(39, 6494, 46022.34375),
    @(40, 6461, 46022.35416666666),
    @(41, 6414, 46022.36458333334),
    @(42, 6386, 46022.375),
    @(43, 6354, 46022.38541666666),
    @(44, 6311, 46022.39583333334)
)

foreach ($row in $data) {
    $r = $row[0]
    $a = $row[1]
    $b = $row[2]
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

